$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 16, pushing the existing rows 16-56 down to 17-57.
$ws.Rows.Item(16).Insert()

# Populate the new row 16 with a new weekly data point (same dimensional
# attributes as the row that used to sit there, new date/price/origin).
$ws.Cells.Item(16, 1).Value = 6
$ws.Cells.Item(16, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(16, 3).Value = "Metropolitana"
$ws.Cells.Item(16, 4).Value = "2021-12-28"
$ws.Cells.Item(16, 5).Value = 13
$ws.Cells.Item(16, 6).Value = "Fruta"
$ws.Cells.Item(16, 7).Value = 100101
$ws.Cells.Item(16, 8).Value = "Berries"
$ws.Cells.Item(16, 9).Value = 100101008
$ws.Cells.Item(16, 10).Value = "Mora"
$ws.Cells.Item(16, 11).Value = "Sin especificar"
$ws.Cells.Item(16, 12).Value = "Primera"
$ws.Cells.Item(16, 13).Value = 100
$ws.Cells.Item(16, 14).Value = 5000
$ws.Cells.Item(16, 15).Value = 6000
$ws.Cells.Item(16, 16).Value = 5500
$ws.Cells.Item(16, 17).Value = "`$/bandeja 2 kilos"
$ws.Cells.Item(16, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(16, 19).Value = 2750
$ws.Cells.Item(16, 20).Value = 2
